# "final draft text & metadata edits before review"
#
# Personnel sheet gets a new collaborator row, an existing collaborator's
# record gets fleshed out with contact details, and the old placeholder
# row for that same person is replaced by a different, new team member.
# A trailing row (existing "Rachel Shrives" record) is appended at the
# bottom so the sheet grows from 9 to 10 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# NOTE ON ORDERING: new text introduces new shared-string-table entries in
# the order cells are first written, so the writes below are sequenced to
# reproduce that same append order (userId, middle initial, given name,
# surname, e-mail) rather than grouped strictly row-by-row.

# ---------------------------------------------------------------------
# 1) Row 9 gains an ORCID-style userId and a middle initial for Zoe
#    Sandwith (given name / surname already present from before).
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 6).Value = "0000-0001-9952-9526"
$ws.Cells.Item(9, 2).Value = "O"

# ---------------------------------------------------------------------
# 2) Row 7 becomes a brand-new collaborator, "S. Alejandra Casillo Cieza"
# ---------------------------------------------------------------------
$ws.Cells.Item(7, 1).Value = "S. Alejandra"
$ws.Cells.Item(7, 3).Value = "Casillo Cieza"
$ws.Cells.Item(7, 4).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(7, 7).Value = "creator"
$ws.Cells.Item(7, 8).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(7, 9).Value = "NSF"
$ws.Cells.Item(7, 10).Value = "OCE-1655686"

# ---------------------------------------------------------------------
# 3) Row 9 also gets a real e-mail, styled as a hyperlink, finishing out
#    Zoe Sandwith's record.
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = "Zoe"
$ws.Cells.Item(9, 3).Value = "Sandwith"
$ws.Cells.Item(9, 4).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(9, 5).Value = "zoe.sandwith@hakai.org"
$ws.Range("E9").Style = "Hyperlink"
$ws.Cells.Item(9, 7).Value = "creator"
$ws.Cells.Item(9, 8).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(9, 9).Value = "NSF"
$ws.Cells.Item(9, 10).Value = "OCE-1655686"

# ---------------------------------------------------------------------
# 4) Row 6 becomes Danielle Aldrett (the old row 7 content); the old
#    row 6 had an empty, hyperlink-styled E cell that needs to go away
#    entirely rather than just lose its value.
# ---------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value = "Danielle"
$ws.Cells.Item(6, 3).Value = "Aldrett"
$ws.Cells.Item(6, 4).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(6, 5).Clear()
$ws.Cells.Item(6, 7).Value = "creator"
$ws.Cells.Item(6, 8).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(6, 9).Value = "NSF"
$ws.Cells.Item(6, 10).Value = "OCE-1655686"

# ---------------------------------------------------------------------
# 5) Append a new row 10 holding what used to be row 9's content
#    (Rachel Shrives / Northeast U.S. Shelf LTER / creator / NSF / OCE-1655686)
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = "Rachel"
$ws.Cells.Item(10, 3).Value = "Shrives"
$ws.Cells.Item(10, 4).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(10, 7).Value = "creator"
$ws.Cells.Item(10, 8).Value = "Northeast U.S. Shelf LTER"
$ws.Cells.Item(10, 9).Value = "NSF"
$ws.Cells.Item(10, 10).Value = "OCE-1655686"

# ---------------------------------------------------------------------
# 6) Leave the final selection on E9, matching where editing wrapped up.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E9").Select()
